$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Numeración de columnas (terminales) en la fila 1: C1=0 ... AL1=35
for ($col = 3; $col -le 38; $col++) {
    $ws.Cells.Item(1, $col).Value = $col - 3
}

# Numeración de filas (no terminales) en la columna A: A3=0 ... A45=42
for ($row = 3; $row -le 45; $row++) {
    $ws.Cells.Item($row, 1).Value = $row - 3
}

# Ensanchar la columna A para que la numeración sea legible
$ws.Columns.Item(1).ColumnWidth = 5.71

# Acercar el zoom de la hoja
$excel.ActiveWindow.Zoom = 40

# Seleccionar el rango de la numeración recién añadida
$ws.Range("A3:A45").Select() | Out-Null
